$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 748.0897
$ws.Range("I17").Value = 979.6
$ws.Range("J17").Value = 714.0441
$ws.Range("K17").Value = 2938.8
$ws.Range("L17").Value = 2142.1323
$ws.Range("M17").Value = -2770.8
$ws.Range("N17").Value = -2478.1323
# Row 19
$ws.Range("H19").Value = 1299.6666
$ws.Range("J19").Value = 1199.6666
$ws.Range("L19").Value = 1199.6666
$ws.Range("N19").Value = -1549.6666
# Row 43
$ws.Range("H43").Value = 4499.8335
$ws.Range("I43").Value = 3399.8
$ws.Range("K43").Value = 3399.8
$ws.Range("M43").Value = -3330.8
# Row 46
$ws.Range("H46").Value = 3192.7144
$ws.Range("I46").Value = 2100
$ws.Range("J46").Value = 3629.8
$ws.Range("K46").Value = 6300
$ws.Range("L46").Value = 10889.4
$ws.Range("M46").Value = -6181
$ws.Range("N46").Value = -11127.4
# Row 55
$ws.Range("H55").Value = 68.333336
$ws.Range("I55").Value = 86.5
$ws.Range("K55").Value = 86.5
$ws.Range("M55").Value = 127.5
# Row 60
$ws.Range("H60").Value = 3192.7144
$ws.Range("I60").Value = 2100
$ws.Range("J60").Value = 3629.8
$ws.Range("K60").Value = 6300
$ws.Range("L60").Value = 10889.4
$ws.Range("M60").Value = -5816
$ws.Range("N60").Value = -11857.4
# Row 64
$ws.Range("H64").Value = 3447.3333
# Row 67
$ws.Range("H67").Value = 3447.3333
# Row 74
$ws.Range("H74").Value = 5599.4443
$ws.Range("I74").Value = 5599.4443
$ws.Range("K74").Value = 5599.4443
$ws.Range("M74").Value = -4663.4443
# Row 77
$ws.Range("H77").Value = 5599.4443
$ws.Range("I77").Value = 5599.4443
$ws.Range("K77").Value = 27997.2215
$ws.Range("M77").Value = -23317.2215
# Row 92
$ws.Range("H92").Value = 2174.2856
$ws.Range("I92").Value = 2174.2856
$ws.Range("K92").Value = 2174.2856
$ws.Range("M92").Value = -926.2856000000002
# Row 107
$ws.Range("H107").Value = 2339.3333
$ws.Range("J107").Value = 718.75
$ws.Range("L107").Value = 718.75
$ws.Range("N107").Value = -4558.75
# Row 132
$ws.Range("H132").Value = 18974.455
$ws.Range("I132").Value = 25641.625
$ws.Range("K132").Value = 76924.875
$ws.Range("M132").Value = -74394.875
# Row 133
$ws.Range("H133").Value = 99990
$ws.Range("J133").Value = 99990
$ws.Range("L133").Value = 99990
$ws.Range("N133").Value = -110110

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1183
$ws.Range("I2").Value = 1178.4
$ws.Range("J2").Value = 1194.5
$ws.Range("K2").Value = 1178.4
$ws.Range("L2").Value = 1194.5
$ws.Range("M2").Value = -1065.4
$ws.Range("N2").Value = -1420.5
# Row 32
$ws.Range("H32").Value = 218530.58
$ws.Range("I32").Value = 221107.33
$ws.Range("K32").Value = 221107.33
$ws.Range("M32").Value = -220820.33
# Row 45
$ws.Range("H45").Value = 3088.9
$ws.Range("I45").Value = 2212.7144
$ws.Range("K45").Value = 2212.7144
$ws.Range("M45").Value = -1835.7144
# Row 63
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2314
# Row 66
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 15000
$ws.Range("M66").Value = -11568
# Row 102
$ws.Range("H102").Value = 240.66667
$ws.Range("I102").Value = 240.66667
$ws.Range("K102").Value = 240.66667
$ws.Range("M102").Value = 1381.33333
# Row 116
$ws.Range("H116").Value = 1183
$ws.Range("I116").Value = 1178.4
$ws.Range("J116").Value = 1194.5
$ws.Range("K116").Value = 1178.4
$ws.Range("L116").Value = 1194.5
$ws.Range("M116").Value = 1115.6
$ws.Range("N116").Value = -5782.5
# Row 139
$ws.Range("H139").Value = 100715
$ws.Range("J139").Value = 100715
$ws.Range("L139").Value = 100715
$ws.Range("N139").Value = -110995

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1183
$ws.Range("I3").Value = 1178.4
$ws.Range("J3").Value = 1194.5
$ws.Range("K3").Value = 1178.4
$ws.Range("L3").Value = 1194.5
$ws.Range("M3").Value = -1064.4
$ws.Range("N3").Value = -1422.5
# Row 99
$ws.Range("H99").Value = 12482.667
$ws.Range("I99").Value = 13893
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 13893
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = -12395
$ws.Range("N99").Value = -4196
# Row 105
$ws.Range("H105").Value = 5237.8184
$ws.Range("I105").Value = 1960.7778
$ws.Range("K105").Value = 1960.7778
$ws.Range("M105").Value = -213.7778000000001
# Row 107
$ws.Range("H107").Value = 1782.5
$ws.Range("I107").Value = 1232
$ws.Range("K107").Value = 1232
$ws.Range("M107").Value = 688
# Row 134
$ws.Range("H134").Value = 6713.857
$ws.Range("I134").Value = 7399.4
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 22198.2
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -19663.2
$ws.Range("N134").Value = -20070

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3841.3
$ws.Range("I31").Value = 3379.2222
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 3379.2222
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -3084.2222
$ws.Range("N31").Value = -8590
# Row 34
$ws.Range("H34").Value = 3841.3
$ws.Range("I34").Value = 3379.2222
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 3379.2222
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -3177.2222
$ws.Range("N34").Value = -8404
# Row 62
$ws.Range("H62").Value = 4002.5
$ws.Range("I62").Value = 4002.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4002.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3378.5
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 4002.5
$ws.Range("I65").Value = 4002.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20012.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -16892.5
$ws.Range("N65").ClearContents()
# Row 92
$ws.Range("H92").Value = 275000
$ws.Range("J92").Value = 275000
$ws.Range("L92").Value = 275000
$ws.Range("N92").Value = -279992
# Row 122
$ws.Range("H122").Value = 12454.48
$ws.Range("I122").Value = 2345.1428
$ws.Range("J122").Value = 65528.5
$ws.Range("K122").Value = 7035.428400000001
$ws.Range("L122").Value = 196585.5
$ws.Range("M122").Value = -4585.428400000001
$ws.Range("N122").Value = -201485.5
# Row 132
$ws.Range("H132").Value = 2060.6365
$ws.Range("I132").Value = 2062.0466
$ws.Range("K132").Value = 6186.139800000001
$ws.Range("M132").Value = -3656.139800000001
# Row 134
$ws.Range("H134").Value = 3268.25
$ws.Range("I134").Value = 2936.6365
$ws.Range("J134").Value = 3997.8
$ws.Range("K134").Value = 8809.9095
$ws.Range("L134").Value = 11993.4
$ws.Range("M134").Value = -6274.9095
$ws.Range("N134").Value = -17063.4

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 29
$ws.Range("I7").Value = 9
$ws.Range("J7").Value = 42.333332
$ws.Range("K7").Value = 27
$ws.Range("L7").Value = 126.999996
$ws.Range("M7").Value = 85
$ws.Range("N7").Value = -350.999996
# Row 113
$ws.Range("H113").Value = 875.5238000000001
$ws.Range("J113").Value = 932.2778
$ws.Range("L113").Value = 2796.8334
$ws.Range("N113").Value = -7136.8334
# Row 127
$ws.Range("H127").Value = 9178.625
$ws.Range("J127").Value = 9178.625
$ws.Range("L127").Value = 27535.875
$ws.Range("N127").Value = -37455.875

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 45000
$ws.Range("J34").Value = 45000
$ws.Range("L34").Value = 45000
$ws.Range("N34").Value = -45536
# Row 76
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45630
# Row 79
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47184
# Row 113
$ws.Range("H113").Value = 2173.3684
$ws.Range("I113").Value = 2183
$ws.Range("K113").Value = 2183
$ws.Range("M113").Value = -13
# Row 122
$ws.Range("H122").Value = 3283.5833
$ws.Range("I122").Value = 3264.5
$ws.Range("J122").Value = 3321.75
$ws.Range("K122").Value = 9793.5
$ws.Range("L122").Value = 9965.25
$ws.Range("M122").Value = -7343.5
$ws.Range("N122").Value = -14865.25

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 46
$ws.Range("H46").Value = 4064.8823
$ws.Range("J46").Value = 4890.7393
$ws.Range("L46").Value = 4890.7393
$ws.Range("N46").Value = -5266.7393

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 44722
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
# Row 107
$ws.Range("H107").Value = 2355.8235
$ws.Range("I107").Value = 1204.381
$ws.Range("K107").Value = 3613.143
$ws.Range("M107").Value = -1693.143
# Row 132
$ws.Range("H132").Value = 2669.577
$ws.Range("I132").Value = 1811.5385
$ws.Range("K132").Value = 5434.6155
$ws.Range("M132").Value = -2904.6155
